# NB train and test fixed (they were being trained with the labels)
# Adds three new worksheets with corrected Naive Bayes evaluation results:
#   ClinicalNB, LabNB, LabNB1
# Each sheet mirrors the existing "CompleteNB" sheet's layout
# (headers f1_score/accuracy in B1/C1, A2 = 0) and carries the
# newly computed f1_score/accuracy pair in B2/C2.

$wb = $excel.ActiveWorkbook
$template = $wb.Worksheets.Item("CompleteNB")

$newSheets = @(
    @{ Name = "ClinicalNB"; F1 = 1; Accuracy = 1 },
    @{ Name = "LabNB";      F1 = 1; Accuracy = 1 },
    @{ Name = "LabNB1";     F1 = 0.4444168646735513; Accuracy = 0.6872 }
)

foreach ($s in $newSheets) {
    $after = $wb.Worksheets.Item($wb.Worksheets.Count)
    $template.Copy($null, $after)
    $newWs = $wb.Worksheets.Item($wb.Worksheets.Count)
    $newWs.Name = $s.Name
    $newWs.Range("B2").Value = $s.F1
    $newWs.Range("C2").Value = $s.Accuracy
}
